{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" problems in the practice\n// table with a new set of problems, cell-by-cell (row, col are the\n// table's own 0-based indices so we never confuse an old value that\n// happens to equal another cell's new value, e.g. \"138\u00f74=\").\nconst replacements = [\n  { row: 0, col: 0, text: \"844\u00f76=\" },\n  { row: 0, col: 1, text: \"865\u00f74=\" },\n  { row: 0, col: 2, text: \"368\u00f73=\" },\n  { row: 0, col: 3, text: \"764\u00f75=\" },\n  { row: 0, col: 4, text: \"799\u00f77=\" },\n\n  { row: 4, col: 0, text: \"156\u00f76=\" },\n  { row: 4, col: 1, text: \"875\u00f74=\" },\n  { row: 4, col: 2, text: \"520\u00f72=\" },\n  { row: 4, col: 3, text: \"198\u00f78=\" },\n  { row: 4, col: 4, text: \"876\u00f75=\" },\n\n  { row: 8, col: 0, text: \"641\u00f74=\" },\n  { row: 8, col: 1, text: \"552\u00f72=\" },\n  { row: 8, col: 2, text: \"191\u00f76=\" },\n  { row: 8, col: 3, text: \"637\u00f74=\" },\n  { row: 8, col: 4, text: \"432\u00f72=\" },\n\n  { row: 12, col: 0, text: \"477\u00f75=\" },\n  { row: 12, col: 1, text: \"402\u00f74=\" },\n  { row: 12, col: 2, text: \"138\u00f74=\" },\n  { row: 12, col: 3, text: \"483\u00f78=\" },\n  { row: 12, col: 4, text: \"315\u00f73=\" },\n\n  { row: 16, col: 0, text: \"796\u00f75=\" },\n  { row: 16, col: 1, text: \"413\u00f78=\" },\n  { row: 16, col: 2, text: \"189\u00f72=\" },\n  { row: 16, col: 3, text: \"888\u00f75=\" },\n  { row: 16, col: 4, text: \"637\u00f75=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const r of replacements) {\n  table.getCell(r.row, r.col).value = r.text;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"three-digit \u00f7 one-digit\" problems in the practice\n# table with a new set of problems. Addressed by (row, col) on the\n# table itself (1-based, COM-style) rather than by text search, so an\n# old value that collides with another cell's new value (e.g.\n# \"138\u00f74=\") can never be double-matched.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n  @{ Row = 1;  Col = 1; Text = \"844\u00f76=\" },\n  @{ Row = 1;  Col = 2; Text = \"865\u00f74=\" },\n  @{ Row = 1;  Col = 3; Text = \"368\u00f73=\" },\n  @{ Row = 1;  Col = 4; Text = \"764\u00f75=\" },\n  @{ Row = 1;  Col = 5; Text = \"799\u00f77=\" },\n\n  @{ Row = 5;  Col = 1; Text = \"156\u00f76=\" },\n  @{ Row = 5;  Col = 2; Text = \"875\u00f74=\" },\n  @{ Row = 5;  Col = 3; Text = \"520\u00f72=\" },\n  @{ Row = 5;  Col = 4; Text = \"198\u00f78=\" },\n  @{ Row = 5;  Col = 5; Text = \"876\u00f75=\" },\n\n  @{ Row = 9;  Col = 1; Text = \"641\u00f74=\" },\n  @{ Row = 9;  Col = 2; Text = \"552\u00f72=\" },\n  @{ Row = 9;  Col = 3; Text = \"191\u00f76=\" },\n  @{ Row = 9;  Col = 4; Text = \"637\u00f74=\" },\n  @{ Row = 9;  Col = 5; Text = \"432\u00f72=\" },\n\n  @{ Row = 13; Col = 1; Text = \"477\u00f75=\" },\n  @{ Row = 13; Col = 2; Text = \"402\u00f74=\" },\n  @{ Row = 13; Col = 3; Text = \"138\u00f74=\" },\n  @{ Row = 13; Col = 4; Text = \"483\u00f78=\" },\n  @{ Row = 13; Col = 5; Text = \"315\u00f73=\" },\n\n  @{ Row = 17; Col = 1; Text = \"796\u00f75=\" },\n  @{ Row = 17; Col = 2; Text = \"413\u00f78=\" },\n  @{ Row = 17; Col = 3; Text = \"189\u00f72=\" },\n  @{ Row = 17; Col = 4; Text = \"888\u00f75=\" },\n  @{ Row = 17; Col = 5; Text = \"637\u00f75=\" }\n)\n\nforeach ($r in $replacements) {\n  $t.Cell($r.Row, $r.Col).Range.Text = $r.Text\n}\n"}
